# Weekly update: a new Chirimoya price record for Vega Monumental Concepción
# is added, and the existing rows 10-25 each shift down by one row
# (row N's data becomes what row N-1 used to hold), with row 10 receiving the
# brand new record and a new row 26 picking up the data that used to be in
# row 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for row 10 (the newly-reported record) ---
$ws.Cells.Item(10, 4).Value = 44804                        # D10 Fecha
$ws.Cells.Item(10, 14).Value = 29000                       # N10 Precio minimo
$ws.Cells.Item(10, 15).Value = 30000                       # O10 Precio maximo
$ws.Cells.Item(10, 16).Value = 29500                       # P10 Precio promedio ponderado
$ws.Cells.Item(10, 17).Value = "$/bandeja 10 kilos"        # Q10 Unidad de comercializacion
$ws.Cells.Item(10, 19).Value = 2950                        # S10 Precio $/Kg
$ws.Cells.Item(10, 20).Value = 10                          # T10 Kg / unidad

# --- Row 11 (was row 10's old values) ---
$ws.Cells.Item(11, 12).Value = "Primera"                   # L11 Calidad
$ws.Cells.Item(11, 14).Value = 2200                        # N11
$ws.Cells.Item(11, 15).Value = 2200                        # O11
$ws.Cells.Item(11, 16).Value = 2200                        # P11
$ws.Cells.Item(11, 19).Value = 2200                        # S11

# --- Row 12 (was row 11's old values) ---
$ws.Cells.Item(12, 4).Value = 44505                        # D12
$ws.Cells.Item(12, 12).Value = "Segunda"                   # L12
$ws.Cells.Item(12, 14).Value = 1800                        # N12
$ws.Cells.Item(12, 15).Value = 1800                        # O12
$ws.Cells.Item(12, 16).Value = 1800                        # P12
$ws.Cells.Item(12, 19).Value = 1800                        # S12

# --- Row 13 (was row 12's old values) ---
$ws.Cells.Item(13, 12).Value = "Primera"                   # L13
$ws.Cells.Item(13, 13).Value = 100                         # M13
$ws.Cells.Item(13, 14).Value = 2700                        # N13
$ws.Cells.Item(13, 15).Value = 2800                        # O13
$ws.Cells.Item(13, 16).Value = 2750                        # P13
$ws.Cells.Item(13, 19).Value = 2750                        # S13

# --- Row 14 (was row 13's old values) ---
$ws.Cells.Item(14, 4).Value = 44467                        # D14
$ws.Cells.Item(14, 12).Value = "Segunda"                   # L14
$ws.Cells.Item(14, 13).Value = 50                          # M14
$ws.Cells.Item(14, 14).Value = 2500                        # N14
$ws.Cells.Item(14, 15).Value = 2500                        # O14
$ws.Cells.Item(14, 16).Value = 2500                        # P14
$ws.Cells.Item(14, 17).Value = "$/kilo (en caja de 15 kilos)"  # Q14
$ws.Cells.Item(14, 19).Value = 2500                        # S14
$ws.Cells.Item(14, 20).Value = 1                           # T14

# --- Row 15 (was row 14's old values) ---
$ws.Cells.Item(15, 12).Value = "Primera"                   # L15
$ws.Cells.Item(15, 13).Value = 100                         # M15
$ws.Cells.Item(15, 14).Value = 17000                       # N15
$ws.Cells.Item(15, 15).Value = 18000                       # O15
$ws.Cells.Item(15, 16).Value = 17500                       # P15
$ws.Cells.Item(15, 19).Value = 2188                        # S15

# --- Row 16 (was row 15's old values) ---
$ws.Cells.Item(16, 4).Value = 44160                        # D16
$ws.Cells.Item(16, 12).Value = "Segunda"                   # L16
$ws.Cells.Item(16, 13).Value = 50                          # M16
$ws.Cells.Item(16, 14).Value = 15000                       # N16
$ws.Cells.Item(16, 15).Value = 15000                       # O16
$ws.Cells.Item(16, 16).Value = 15000                       # P16
$ws.Cells.Item(16, 17).Value = "$/bandeja 8 kilos"         # Q16
$ws.Cells.Item(16, 19).Value = 1875                        # S16
$ws.Cells.Item(16, 20).Value = 8                           # T16

# --- Row 17 (was row 16's old values) ---
$ws.Cells.Item(17, 4).Value = 44517                        # D17
$ws.Cells.Item(17, 14).Value = 25000                       # N17
$ws.Cells.Item(17, 15).Value = 27000                       # O17
$ws.Cells.Item(17, 16).Value = 26000                       # P17
$ws.Cells.Item(17, 19).Value = 2600                        # S17

# --- Row 18 (was row 17's old values) ---
$ws.Cells.Item(18, 4).Value = 44461                        # D18
$ws.Cells.Item(18, 14).Value = 29000                       # N18
$ws.Cells.Item(18, 15).Value = 30000                       # O18
$ws.Cells.Item(18, 16).Value = 29500                       # P18
$ws.Cells.Item(18, 19).Value = 2950                        # S18

# --- Row 19 (was row 18's old values) ---
$ws.Cells.Item(19, 4).Value = 44469                        # D19
$ws.Cells.Item(19, 13).Value = 100                         # M19
$ws.Cells.Item(19, 14).Value = 28000                       # N19
$ws.Cells.Item(19, 15).Value = 29000                       # O19
$ws.Cells.Item(19, 16).Value = 28500                       # P19
$ws.Cells.Item(19, 19).Value = 2850                        # S19

# --- Row 20 (was row 19's old values) ---
$ws.Cells.Item(20, 4).Value = 44488                        # D20
$ws.Cells.Item(20, 13).Value = 50                          # M20
$ws.Cells.Item(20, 14).Value = 25000                       # N20
$ws.Cells.Item(20, 15).Value = 26000                       # O20
$ws.Cells.Item(20, 16).Value = 25600                       # P20
$ws.Cells.Item(20, 17).Value = "$/bandeja 10 kilos"        # Q20
$ws.Cells.Item(20, 19).Value = 2560                        # S20
$ws.Cells.Item(20, 20).Value = 10                          # T20

# --- Row 21 (was row 20's old values) ---
$ws.Cells.Item(21, 12).Value = "Primera"                   # L21
$ws.Cells.Item(21, 13).Value = 100                         # M21
$ws.Cells.Item(21, 14).Value = 2900                        # N21
$ws.Cells.Item(21, 15).Value = 3000                        # O21
$ws.Cells.Item(21, 16).Value = 2950                        # P21
$ws.Cells.Item(21, 19).Value = 2950                        # S21

# --- Row 22 (was row 21's old values) ---
$ws.Cells.Item(22, 4).Value = 44462                        # D22
$ws.Cells.Item(22, 12).Value = "Segunda"                   # L22
$ws.Cells.Item(22, 13).Value = 50                          # M22
$ws.Cells.Item(22, 14).Value = 2600                        # N22
$ws.Cells.Item(22, 15).Value = 2600                        # O22
$ws.Cells.Item(22, 16).Value = 2600                        # P22
$ws.Cells.Item(22, 17).Value = "$/kilo (en caja de 15 kilos)"  # Q22
$ws.Cells.Item(22, 19).Value = 2600                        # S22
$ws.Cells.Item(22, 20).Value = 1                           # T22

# --- Row 23 (was row 22's old values) ---
$ws.Cells.Item(23, 4).Value = 44454                        # D23
$ws.Cells.Item(23, 14).Value = 30000                       # N23
$ws.Cells.Item(23, 15).Value = 31000                       # O23
$ws.Cells.Item(23, 16).Value = 30500                       # P23
$ws.Cells.Item(23, 19).Value = 3050                        # S23

# --- Row 24 (was row 23's old values) ---
$ws.Cells.Item(24, 4).Value = 44484                        # D24
$ws.Cells.Item(24, 14).Value = 25000                       # N24
$ws.Cells.Item(24, 15).Value = 26000                       # O24
$ws.Cells.Item(24, 16).Value = 25500                       # P24
$ws.Cells.Item(24, 17).Value = "$/bandeja 10 kilos"        # Q24
$ws.Cells.Item(24, 19).Value = 2550                        # S24
$ws.Cells.Item(24, 20).Value = 10                          # T24

# --- Row 25 (was row 24's old values) ---
$ws.Cells.Item(25, 12).Value = "Primera"                   # L25
$ws.Cells.Item(25, 13).Value = 100                         # M25
$ws.Cells.Item(25, 14).Value = 1900                        # N25
$ws.Cells.Item(25, 15).Value = 2000                        # O25
$ws.Cells.Item(25, 16).Value = 1950                        # P25
$ws.Cells.Item(25, 19).Value = 1950                        # S25

# --- New row 26 (was row 25's old values, now appended as a new row) ---
$ws.Cells.Item(26, 1).Value = 11                                   # A26 Mercado ID
$ws.Cells.Item(26, 2).Value = "Vega Monumental Concepción"         # B26 Mercado
$ws.Cells.Item(26, 3).Value = "Bíobío"                             # C26 Región
$ws.Cells.Item(26, 4).Value = 44516                                # D26 Fecha
$ws.Cells.Item(26, 5).Value = 8                                    # E26 Codreg
$ws.Cells.Item(26, 6).Value = "Fruta"                              # F26 Tipo
$ws.Cells.Item(26, 7).Value = 100107                               # G26 Producto ID
$ws.Cells.Item(26, 8).Value = "Otros"                              # H26 Producto
$ws.Cells.Item(26, 9).Value = 100107002                            # I26 Categoria ID
$ws.Cells.Item(26, 10).Value = "Chirimoya"                         # J26 Categoria
$ws.Cells.Item(26, 11).Value = "Cultivar IV Región"                # K26 Variedad
$ws.Cells.Item(26, 12).Value = "Segunda"                           # L26 Calidad
$ws.Cells.Item(26, 13).Value = 50                                  # M26 Volumen
$ws.Cells.Item(26, 14).Value = 1700                                # N26 Precio minimo
$ws.Cells.Item(26, 15).Value = 1700                                # O26 Precio maximo
$ws.Cells.Item(26, 16).Value = 1700                                # P26 Precio promedio ponderado
$ws.Cells.Item(26, 17).Value = "$/kilo (en caja de 15 kilos)"      # Q26 Unidad de comercializacion
$ws.Cells.Item(26, 18).Value = "Provincia de Limarí"               # R26 Origen
$ws.Cells.Item(26, 19).Value = 1700                                # S26 Precio $/Kg
$ws.Cells.Item(26, 20).Value = 1                                   # T26 Kg / unidad

# Also make sure D26's date style matches the other date cells in column D
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat
